# Insert two new data rows at row 74 (pushing the existing rows 74-114 down
# to 76-116), then populate the two new rows with the new weekly records.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(74).Resize(2, 1).EntireRow.Insert()

# New row 74
$ws.Cells.Item(74, 1).Value = 10
$ws.Cells.Item(74, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(74, 3).Value = "La Araucanía"
$ws.Cells.Item(74, 4).Value = 44609
$ws.Cells.Item(74, 5).Value = 9
$ws.Cells.Item(74, 6).Value = 100112031
$ws.Cells.Item(74, 7).Value = "Poroto verde"
$ws.Cells.Item(74, 8).Value = "Brío"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 80
$ws.Cells.Item(74, 11).Value = 1200
$ws.Cells.Item(74, 12).Value = 1200
$ws.Cells.Item(74, 13).Value = 1200
$ws.Cells.Item(74, 14).Value = "`$/kilo"
$ws.Cells.Item(74, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(74, 16).Value = 1200
$ws.Cells.Item(74, 17).Value = 1
$ws.Cells.Item(74, 18).Value = "Hortaliza"

# New row 75
$ws.Cells.Item(75, 1).Value = 10
$ws.Cells.Item(75, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(75, 3).Value = "La Araucanía"
$ws.Cells.Item(75, 4).Value = 44609
$ws.Cells.Item(75, 5).Value = 9
$ws.Cells.Item(75, 6).Value = 100112031
$ws.Cells.Item(75, 7).Value = "Poroto verde"
$ws.Cells.Item(75, 8).Value = "Sin especificar"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 80
$ws.Cells.Item(75, 11).Value = 1200
$ws.Cells.Item(75, 12).Value = 1200
$ws.Cells.Item(75, 13).Value = 1200
$ws.Cells.Item(75, 14).Value = "`$/kilo"
$ws.Cells.Item(75, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(75, 16).Value = 1200
$ws.Cells.Item(75, 17).Value = 1
$ws.Cells.Item(75, 18).Value = "Hortaliza"
